# Rows 63-65 hold three separate sighting records that got re-sorted:
# each row now carries the data that used to belong to the row above it,
# with row 63 wrapping around to what used to be on row 65 (i.e. a cyclic
# shift 63 -> 64 -> 65 -> 63 of the record payloads, the row-independent
# columns such as lan/kommun/rapportor stay where they are).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column "I" ("Antal") is stored as text even when its content looks like
# a plain integer (e.g. "3"), so force text formatting on it before
# writing back - otherwise a bare numeric assignment would turn it into
# a genuine number cell.
$ws.Range("I63:I65").NumberFormat = "@"

# --- snapshot current (pre-edit) values for rows 63, 64, 65 -----------
function Get-RowData($r) {
    return @{
        A = $ws.Range("A$r").Value()
        B = $ws.Range("B$r").Value()
        E = $ws.Range("E$r").Value()
        F = $ws.Range("F$r").Value()
        G = $ws.Range("G$r").Value()
        H = $ws.Range("H$r").Value()
        I = $ws.Range("I$r").Value()
        J = $ws.Range("J$r").Value()
        P = $ws.Range("P$r").Value()
        Q = $ws.Range("Q$r").Value()
        R = $ws.Range("R$r").Value()
        S = $ws.Range("S$r").Value()
        Z = $ws.Range("Z$r").Value()
        AB = $ws.Range("AB$r").Value()
    }
}

$row63 = Get-RowData 63
$row64 = Get-RowData 64
$row65 = Get-RowData 65

function Set-RowData($r, $data) {
    $ws.Range("A$r").Value = $data.A
    $ws.Range("B$r").Value = $data.B
    $ws.Range("E$r").Value = $data.E
    $ws.Range("F$r").Value = $data.F
    $ws.Range("G$r").Value = $data.G
    $ws.Range("H$r").Value = $data.H
    $ws.Range("I$r").Value = $data.I
    $ws.Range("J$r").Value = $data.J
    $ws.Range("P$r").Value = $data.P
    $ws.Range("Q$r").Value = $data.Q
    $ws.Range("R$r").Value = $data.R
    $ws.Range("S$r").Value = $data.S
    $ws.Range("Z$r").Value = $data.Z
    $ws.Range("AB$r").Value = $data.AB
}

# row 63 <- old row 65, row 64 <- old row 63, row 65 <- old row 64
Set-RowData 63 $row65
Set-RowData 64 $row63
Set-RowData 65 $row64

# The "Kon" (L) value that used to sit on row 64 (blank) moves down to
# row 65; row 64 no longer carries that column at all.
$ws.Range("L64").ClearContents()
